$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CCI_INSPQ_2018_CIM9" — selection only changes (C1 -> A32)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("CCI_INSPQ_2018_CIM9")
$ws1.Range("A32").Select()

# ---------------------------------------------------------------------------
# Sheet "CCI_INSPQ_2018_CIM10" — selection only changes (C1 -> A32:C32)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CCI_INSPQ_2018_CIM10")
$ws2.Range("A32:C32").Select()

# ---------------------------------------------------------------------------
# Sheet "UManitoba_2016" — re-sort the data (A2:C18) by the abbreviation
# column (B) instead of the diagnosis-label column (A); no data changes.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("UManitoba_2016")
$sort3 = $ws3.Sort
$sort3.SortFields.Clear()
$sort3.SortFields.Add($ws3.Range("B2:B18"))
$sort3.SetRange($ws3.Range("A1:C18"))
$sort3.Header = 1
$sort3.Apply()
$ws3.Range("B38").Select()

# ---------------------------------------------------------------------------
# New sheet "CCI_INSPQ_Manitoba" — added after UManitoba_2016, combining the
# Charlson/INSPQ abbreviations with Manitoba-style weights.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "CCI_INSPQ_Manitoba"

# Header row (same three labels used on every other sheet)
$ws4.Range("A1").Value = "DIAGN"
$ws4.Range("B1").Value = "DIAGN_CODE"
$ws4.Range("C1").Value = "POIDS"
$ws4.Range("A1:C1").Font.Bold = $true
$ws4.Range("A1:C1").HorizontalAlignment = -4108

# Data rows
$data = @(
    @("HIV/AIDS", "aids", 6),
    @("Cancer", "canc", 2),
    @("Cerebrovascular disease", "cevd", 1),
    @("Congestive Heart Failure", "chf", 1),
    @("Chronic pulmonary disease", "copd", 1),
    @("Dementia", "dementia", 1),
    @("Diabetes with Chronic Complications", "diab", 2),
    @("Diabetes without Chronic Complications", "diabwc", 1),
    @("Mild Liver Disease", "ld1", 1),
    @("Moderate or Severe Liver Disease", "ld2", 3),
    @("Metastatic cancer", "metacanc", 6),
    @("Myocardial Infarction", "mi", 1),
    @("Paralysis", "para", 2),
    @("Peripheral Vascular Disease", "pvd", 1),
    @("Rheumatoid arth./collagen vascular disease", "rheumd", 1),
    @("Ulcer disease", "ud", 1),
    @("Valvular disease", "valv", 0)
)

$r = 2
foreach ($row in $data) {
    $ws4.Cells.Item($r, 1).Value = $row[0]
    $ws4.Cells.Item($r, 2).Value = $row[1]
    $ws4.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Last row (Valvular disease / 0) is formatted as an integer, like the
# corresponding rows on the CIM9/CIM10 sheets.
$ws4.Range("C18").NumberFormat = "0"

# Column widths to match the sister sheets
$ws4.Columns("A").ColumnWidth = $ws3.Columns("A").ColumnWidth
$ws4.Columns("B").ColumnWidth = $ws3.Columns("B").ColumnWidth
$ws4.Columns("C").ColumnWidth = $ws3.Columns("C").ColumnWidth

$ws4.PageSetup.Orientation = 1

# Sort A2:C17 (not including the trailing Valvular-disease row) by the
# abbreviation column (B), matching the other sheets' sort convention.
$sort4 = $ws4.Sort
$sort4.SortFields.Clear()
$sort4.SortFields.Add($ws4.Range("B2:B17"))
$sort4.SetRange($ws4.Range("A2:C17"))
$sort4.Header = 0
$sort4.Apply()

$ws4.Range("C17").Select()

Write-Output "done"
